$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) CAPSTONE_004 + "V" -> merge into a single run "CAPSTONE_004V"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("CAPSTONE_004V", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "CAPSTONE_004V", 2) | Out-Null

# ---------------------------------------------------------------------
# Helper: replace a RUT string with a new RUT, materialised as one run
# per token (digits / dots / dash) the way real Word leaves behind
# separately-typed runs that happen to share identical formatting.
#
# We do this by momentarily nudging each token's Font.Size to a distinct
# value (forcing the engine to split the run at that boundary) and then
# restoring every token to the original size (10pt / sz=20) - the split
# points stay in place even though the final formatting is identical
# across all the resulting runs, matching the target XML.
# ---------------------------------------------------------------------

function Replace-Rut($oldRut, $tokens) {
    $rng = $d.Content
    $found = $rng.Find.Execute($oldRut, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if (-not $found) {
        return
    }

    $newText = [string]::Join("", $tokens)
    $start = $rng.Start
    $rng.Text = $newText

    # Compute the [start,end) offsets of every token within the range.
    $offsets = @()
    $pos = $start
    foreach ($tok in $tokens) {
        $offsets += , @($pos, $pos + $tok.Length)
        $pos = $pos + $tok.Length
    }

    # Force a run boundary at every token by giving each a distinct size.
    $bump = 11
    foreach ($off in $offsets) {
        $sub = $d.Range($off[0], $off[1])
        $sub.Font.Size = $bump
        $bump = $bump + 1
    }

    # Restore the original size (10pt == w:sz 20) on every token; the
    # run split already happened above and persists.
    foreach ($off in $offsets) {
        $sub = $d.Range($off[0], $off[1])
        $sub.Font.Size = 10
    }
}

# ---------------------------------------------------------------------
# 2) RUT changes
# ---------------------------------------------------------------------
Replace-Rut "20.298.598-K" @("12", ".", "900", ".", "628", "-", "5")
Replace-Rut "18.056.442-K" @("17", ".", "812", ".", "454", "-", "4")
